$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.693.15"
$ws.Range("E2").Value = "  -2.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.809.43"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.06"
$ws.Range("E5").Value = "  -2.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.53"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.809.33"
$ws.Range("E7").Value = "  +1.86%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  -3.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.20"
$ws.Range("E11").Value = "  -4.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.464"
$ws.Range("E12").Value = "  -4.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.33"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000245"
$ws.Range("E14").Value = "  -3.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.445.28"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.810.55"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.666.68"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.21"
$ws.Range("E18").Value = "  -3.27%  "
$ws.Range("E19").Value = "  -4.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.60"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "492.23"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.11"
$ws.Range("E22").Value = "  -2.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.745"
$ws.Range("E23").Value = "  +3.42%  "
$ws.Range("E24").Value = "  +14.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.03"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("E26").Value = "  -6.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.30"
$ws.Range("E27").Value = "  -4.35%  "
$ws.Range("E28").Value = "  -6.02%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +0.90%  "
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.67"
$ws.Range("E32").Value = "  +7.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.65"
$ws.Range("E33").Value = "  -4.21%  "
$ws.Range("E34").Value = "  -3.92%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("E36").Value = "  -3.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.80"
$ws.Range("E37").Value = "  -5.10%  "
$ws.Range("E38").Value = "  -2.21%  "
$ws.Range("E39").Value = "  -4.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "451.71"
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "49.17"
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.90"
$ws.Range("E43").Value = "  -3.66%  "
$ws.Range("E44").Value = "  -2.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.30"
$ws.Range("E45").Value = "  -7.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.846.94"
$ws.Range("E46").Value = "  -3.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "139.50"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0350"
$ws.Range("E49").Value = "  -2.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.55"
$ws.Range("E50").Value = "  +13.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.01"
$ws.Range("E51").Value = "  -4.98%  "
